$wb = $excel.ActiveWorkbook

# ---- Sheet2: add weather data rows ----
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Cells.Item(2,1).Value = "Auckland Central"
$ws2.Cells.Item(2,2).Value = "Mon 9 Oct"
$ws2.Cells.Item(2,3).Value = "16°"
$ws2.Cells.Item(2,4).Value = "12°"

$ws2.Cells.Item(3,1).Value = "Hunua"
$ws2.Cells.Item(3,2).Value = "Mon 9 Oct"
$ws2.Cells.Item(3,3).Value = "17°"
$ws2.Cells.Item(3,4).Value = "11°"

$ws2.Cells.Item(4,1).Value = "Kumeu"
$ws2.Cells.Item(4,2).Value = "Mon 9 Oct"
$ws2.Cells.Item(4,3).Value = "16°"
$ws2.Cells.Item(4,4).Value = "11°"

$ws2.Cells.Item(5,1).Value = "Eastern Rangitaiki"
$ws2.Cells.Item(5,2).Value = "Mon 9 Oct"
$ws2.Cells.Item(5,3).Value = "14°"
$ws2.Cells.Item(5,4).Value = "4°"

$ws2.Cells.Item(6,1).Value = "Methven"
$ws2.Cells.Item(6,2).Value = "Mon 9 Oct"
$ws2.Cells.Item(6,3).Value = "17°"
$ws2.Cells.Item(6,4).Value = "4°"

# column B a bit wider, like the author did (engine quantizes column width to
# 1/6-character steps, so 13.6 is the closest we can land to the recorded 14.43)
$ws2.Columns.Item(2).ColumnWidth = 13.6

# ---- Sheet3: drop the personal email, rename header, keep group email ----
$ws3 = $wb.Worksheets.Item("Sheet3")

# clear out hyperlinks before we shuffle rows around (avoids stale refs)
$ws3.Range("A1:A3").Hyperlinks.Delete()

# remove the samft223@gmail.com row entirely, shifting SOFTENG762... up
$ws3.Rows.Item(2).Delete()

$ws3.Range("A1").Value = "emails"
$ws3.Range("A1").Font.Bold = $true

$ws3.Hyperlinks.Add($ws3.Range("A2"), "mailto:SOFTENG762Group10@gmail.com")

# ---- selection / active-tab bookkeeping ----
[void]$ws2.Activate()
[void]$ws2.Range("D3").Select()

[void]$ws3.Activate()
[void]$ws3.Range("B8").Select()
